$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) from Excel auto-converting numeric-looking
# text (e.g. "207.71") into actual numbers -- the source data stores
# these as plain text (inline strings), so force Text format while
# writing, then restore the default "Normal" style so no extra
# formatting is left behind on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.041.34'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '1.560.95'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("D5").Value = '207.71'
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '22.09'
$ws.Range("E8").Value = '  +3.50%  '
$ws.Range("D9").Value = '0.248'
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = '0.0587'
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").Value = '1.786.83'
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").Value = '1.564.39'
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("E14").Value = '  +2.08%  '
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '27.073.87'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '219.31'
$ws.Range("E18").Value = '  +2.26%  '
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").Value = '7.34'
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  +1.62%  '
$ws.Range("D23").Value = '9.26'
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").Value = '154.51'
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("D26").Value = '6.61'
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").Value = '14.97'
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = '1.446.25'
$ws.Range("E33").Value = '  +5.81%  '
$ws.Range("D34").Value = '3.08'
$ws.Range("E34").Value = '  +4.92%  '
$ws.Range("E35").Value = '  +3.88%  '
$ws.Range("D36").Value = '0.963'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  +0.86%  '
$ws.Range("D39").Value = '0.524'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.72'
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("E43").Value = '  +3.36%  '
$ws.Range("D44").Value = '0.988'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("D45").Value = '64.26'
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("D47").Value = '1.698.22'
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("D48").Value = '86.74'
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("D49").Value = '0.0525'
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("E51").Value = '  +2.47%  '

$ws.Range("D2:D51").Style = "Normal"
